$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.919.74"

$ws.Range("D3").Value = "1.810.87"
$ws.Range("E3").Value = "  +2.92%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'311.86"
$ws.Range("E5").Value = "  -2.20%  "

$ws.Range("E6").Value = "  +0.19%  "

$ws.Range("D7").Value = "'0.4288"
$ws.Range("E7").Value = "  -1.08%  "

$ws.Range("D8").Value = "'0.3692"
$ws.Range("E8").Value = "  +2.01%  "

$ws.Range("D9").Value = "'0.07239"
$ws.Range("E9").Value = "  +2.08%  "

$ws.Range("D10").Value = "'0.8620"
$ws.Range("E10").Value = "  +3.81%  "

$ws.Range("D11").Value = "2.043.93"
$ws.Range("E11").Value = "  +15.47%  "

$ws.Range("D12").Value = "'21.14"
$ws.Range("E12").Value = "  +5.14%  "

$ws.Range("D13").Value = "'6.639"
$ws.Range("E13").Value = "  +5.23%  "

$ws.Range("D14").Value = "'5.394"
$ws.Range("E14").Value = "  +3.64%  "

$ws.Range("D15").Value = "'0.06919"
$ws.Range("E15").Value = "  +1.38%  "

$ws.Range("D16").Value = "'80.80"
$ws.Range("E16").Value = "  +2.51%  "

$ws.Range("E17").Value = "  +0.21%  "

$ws.Range("D18").Value = "'0.000008927"
$ws.Range("E18").Value = "  +3.51%  "

$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("E20").Value = "  +2.24%  "

$ws.Range("D21").Value = "26.967.39"
$ws.Range("E21").Value = "  +1.38%  "

$ws.Range("D22").Value = "'5.199"
$ws.Range("E22").Value = "  +4.34%  "

$ws.Range("D23").Value = "'10.95"
$ws.Range("E23").Value = "  -1.15%  "

$ws.Range("D24").Value = "2.272.05"
$ws.Range("E24").Value = "  +14.10%  "

$ws.Range("D25").Value = "'154.10"
$ws.Range("E25").Value = "  +0.95%  "

$ws.Range("D26").Value = "'1.887"
$ws.Range("E26").Value = "  -0.73%  "

$ws.Range("D27").Value = "'18.33"
$ws.Range("E27").Value = "  +1.32%  "

$ws.Range("D28").Value = "'5.234"
$ws.Range("E28").Value = "  +4.16%  "

$ws.Range("D29").Value = "'1.906"
$ws.Range("E29").Value = "  +16.57%  "

$ws.Range("D30").Value = "'114.95"
$ws.Range("E30").Value = "  +0.86%  "

$ws.Range("D31").Value = "'0.08921"
$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("D32").Value = "'0.7418"
$ws.Range("E32").Value = "  +4.07%  "

$ws.Range("D33").Value = "'1.154"
$ws.Range("E33").Value = "  +5.59%  "

$ws.Range("D34").Value = "'4.420"
$ws.Range("E34").Value = "  +2.93%  "

$ws.Range("D35").Value = "'2.807"
$ws.Range("E35").Value = "  +1.31%  "

$ws.Range("D36").Value = "'1.008"
$ws.Range("E36").Value = "  +0.31%  "

$ws.Range("D37").Value = "'1.124"
$ws.Range("E37").Value = "  +5.50%  "

$ws.Range("E38").Value = "  +3.01%  "

$ws.Range("D39").Value = "'0.01924"
$ws.Range("E39").Value = "  +1.60%  "

$ws.Range("D40").Value = "'0.5084"
$ws.Range("E40").Value = "  +4.14%  "

$ws.Range("D41").Value = "'2.757"
$ws.Range("E41").Value = "  +9.72%  "

$ws.Range("D42").Value = "'0.1645"
$ws.Range("E42").Value = "  +3.34%  "

$ws.Range("D43").Value = "'6.436"
$ws.Range("E43").Value = "  +4.82%  "

$ws.Range("D44").Value = "'8.256"
$ws.Range("E44").Value = "  +4.97%  "

$ws.Range("D45").Value = "'107.03"
$ws.Range("E45").Value = "  +2.59%  "

$ws.Range("E46").Value = "  +3.60%  "

$ws.Range("D47").Value = "'1.008"
$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.4575"
$ws.Range("E48").Value = "  +2.80%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.648"
$ws.Range("E49").Value = "  +5.57%  "

$ws.Range("D50").Value = "'0.06282"
$ws.Range("E50").Value = "  +1.49%  "

$ws.Range("E51").Value = "  +5.15%  "

